$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.00" or
# "68.598.03" are not auto-converted to numbers by Excel's smart-typing.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.598.03"
$ws.Range("E2").Value = "  +1.14%  "

$ws.Range("D3").Value = "3.822.61"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "612.15"
$ws.Range("E5").Value = "  +1.27%  "

$ws.Range("D6").Value = "164.57"
$ws.Range("E6").Value = "  -0.95%  "

$ws.Range("D7").Value = "3.821.27"
$ws.Range("E7").Value = "  +0.30%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").Value = "0.518"
$ws.Range("E9").Value = "  -0.22%  "

$ws.Range("E10").Value = "  +0.26%  "

$ws.Range("D11").Value = "0.452"
$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").Value = "6.80"
$ws.Range("E12").Value = "  +7.40%  "

$ws.Range("E13").Value = "  -1.41%  "

$ws.Range("D14").Value = "35.37"
$ws.Range("E14").Value = "  -1.89%  "

$ws.Range("D15").Value = "4.461.49"
$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("D16").Value = "3.877.38"
$ws.Range("E16").Value = "  +1.54%  "

$ws.Range("D17").Value = "68.571.63"
$ws.Range("E17").Value = "  +1.06%  "

$ws.Range("D18").Value = "18.19"
$ws.Range("E18").Value = "  -1.14%  "

$ws.Range("D19").Value = "7.14"
$ws.Range("E19").Value = "  +0.49%  "

$ws.Range("E20").Value = "  -0.10%  "

$ws.Range("D21").Value = "464.06"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").Value = "9.66"
$ws.Range("E22").Value = "  -2.04%  "

$ws.Range("D23").Value = "0.703"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").Value = "0.0000148"
$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "83.70"
$ws.Range("E25").Value = "  +0.33%  "

$ws.Range("D26").Value = "12.04"
$ws.Range("E26").Value = "  -0.97%  "

$ws.Range("D27").Value = "2.12"
$ws.Range("E27").Value = "  -0.38%  "

$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("D29").Value = "10.01"
$ws.Range("E29").Value = "  -0.34%  "

$ws.Range("D30").Value = "3.963.40"
$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("E31").Value = "  -5.65%  "

$ws.Range("D32").Value = "2.22"
$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("D33").Value = "7.26"
$ws.Range("E33").Value = "  -2.04%  "

$ws.Range("D34").Value = "29.01"
$ws.Range("E34").Value = "  -1.81%  "

$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("E37").Value = "  +1.14%  "

$ws.Range("D38").Value = "0.147"
$ws.Range("E38").Value = "  +6.53%  "

$ws.Range("D39").Value = "5.90"
$ws.Range("E39").Value = "  +1.42%  "

$ws.Range("E40").Value = "  -1.49%  "

$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").Value = "3.14"
$ws.Range("E42").Value = "  -3.15%  "

$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("D44").Value = "154.05"
$ws.Range("E44").Value = "  +1.56%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.298"
$ws.Range("E45").Value = "  -0.77%  "

$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").Value = "43.30"
$ws.Range("E46").Value = "  -3.51%  "

$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "1.40"
$ws.Range("E47").Value = "  +0.33%  "

$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "46.61"
$ws.Range("E48").Value = "  -2.37%  "

$ws.Range("D49").Value = "8.39"
$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("D50").Value = "1.87"
$ws.Range("E50").Value = "  +0.88%  "

$ws.Range("D51").Value = "379.55"
$ws.Range("E51").Value = "  -2.93%  "

# Restore default "Normal" style on the Price column so we don't leave
# behind a stray number-format style that wasn't in the original workbook.
$ws.Range("D2:D51").Style = "Normal"
